$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Past Updates" row gains two explanatory notes (E12/F12) ---
$ws.Range("E12").Value = "All past updates should be linked to this page"
$ws.Range("F12").Value = "All internal links in past update pages need to be checked (and some added; older updates don't have links to the published works like new updates do). Formats are also not standardized; we should consider whether we want to do this."

# --- Fill in now-available bibliography / articles links ---
$ws.Range("D23").Value = "generalbib.html"
$ws.Range("D24").Value = "specificbib.html"
$ws.Range("D35").Value = "articlesindex.html"

# --- Insert a new submenu row for "Recent and Forthcoming Presentations" ---
# (pushes "Plan of the Archive" and everything below down by one row)
$ws.Rows("36:36").Insert()
$ws.Range("C36").Value = "Recent and Forthcoming Presentations"
$ws.Range("D36").Value = "presentations.html"

# --- Small wording fix: move the closing quote after "IATH/CDLA" ---
$ws.Range("F38").Value = "This would replace the ""Blake Archive in the Context of IATH/CDLA"" page. It may not be necessary."

# --- Widen column F slightly to fit the newly added, longer notes ---
$ws.Columns("F").ColumnWidth = 56.5
